$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-21 Friday", "2024-06-22 Saturday"),
    @("54×51=2754", "32×82=2624"),
    @("70×88=6160", "89×82=7298"),
    @("41×58=2378", "24×65=1560"),
    @("59×36=2124", "40×77=3080"),
    @("61×29=1769", "87×38=3306"),
    @("73×72=5256", "23×55=1265"),
    @("81×28=2268", "40×24=960"),
    @("72×67=4824", "84×36=3024"),
    @("25×44=1100", "49×62=3038"),
    @("24×87=2088", "34×61=2074"),
    @("94×94=8836", "61×84=5124"),
    @("29×26=754", "46×43=1978"),
    @("22×61=1342", "18×74=1332"),
    @("38×98=3724", "57×42=2394"),
    @("71×52=3692", "39×87=3393"),
    @("23×40=920", "89×88=7832"),
    @("93×68=6324", "72×91=6552"),
    @("58×69=4002", "23×71=1633"),
    @("50×31=1550", "33×78=2574"),
    @("12×43=516", "49×16=784"),
    @("89×41=3649", "67×43=2881"),
    @("13×76=988", "37×57=2109"),
    @("45×19=855", "60×35=2100"),
    @("36×80=2880", "77×24=1848"),
    @("90×91=8190", "28×35=980")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
